$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.495.95"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "3.065.98"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'548.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "'140.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.13%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.059.04"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").Value = "'6.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").Value = "'34.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "3.561.88"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "63.373.21"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "3.068.25"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "'483.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("D21").Value = "'13.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  +4.76%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "'12.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  +5.93%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'2.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.36%  "
$ws.Range("D34").Value = "'5.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("D35").Value = "'55.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "'466.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("D39").Value = "'0.0396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.19%  "
$ws.Range("D40").Value = "3.067.06"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").Value = "'2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").Value = "'28.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").Value = "'116.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "0.0₃0508"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("E51").Value = "  +3.50%  "
